$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46, shifting existing rows 46:143 down to 47:144
$ws.Rows("46:46").Insert()

# Populate the newly inserted row 46 with its data
$ws.Range("A46").Value = 11
$ws.Range("B46").Value = "Vega Monumental Concepción"
$ws.Range("C46").Value = "Bíobío"
$ws.Range("D46").Value = 44614
$ws.Range("E46").Value = 8
$ws.Range("F46").Value = 100112003
$ws.Range("G46").Value = "Ajo"
$ws.Range("H46").Value = "Chino"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 400
$ws.Range("K46").Value = 17000
$ws.Range("L46").Value = 18000
$ws.Range("M46").Value = 17500
$ws.Range("N46").Value = "$/caja 10 kilos"
$ws.Range("O46").Value = "China"
$ws.Range("P46").Value = 1750
$ws.Range("Q46").Value = 10
$ws.Range("R46").Value = "Hortaliza"
